# "fixed viper sting bug" - add games 1..6 counters to the COMPS log on
# "Лист1" and fill in the WIN/LOSE results that were previously missing
# (or mis-recorded) for each game, per the viper-sting bug fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

# Game 1 (rows 19-22): Beast Hunter / Demon Lock / Retri Paladin -> LOSE 2
$ws.Range("A19").Value = 1
$ws.Range("C20").Value = "LOSE 2"

# Game 2 (rows 23-26): Disc Priest / Affli Lock / Frost Mage
$ws.Range("A23").Value = 2
$ws.Range("C24").Value = "WIN 3"
$ws.Range("D24").Value = "WIN 5"
$ws.Range("C25").Value = "WIN 1"
$ws.Range("D25").Value = "LOSE 6"
$ws.Range("C26").Value = "WIN 4"
$ws.Range("J26").Value = "skipping at max mana gives you more mana than the max…"

# Game 3 (rows 27-30): Surv Hunter / Disc Priest / Holy Paladin
$ws.Range("A27").Value = 3
$ws.Range("C28").Value = "LOSE 2"
$ws.Range("C29").Value = "LOSE 6"

# Game 4 (rows 31-34): Sub Rogue / Frost Mage / Retri Paladin -> LOSE 2
$ws.Range("A31").Value = 4
$ws.Range("C32").Value = "LOSE 2"

# Game 5 (rows 35-38): Disc Priest / Frost Mage / Retri Paladin -> LOSE 3
$ws.Range("A35").Value = 5
$ws.Range("A36").Value = "DISC PRIEST"
$ws.Range("C36").Value = "LOSE 3"
$ws.Range("A37").Value = "FROST MAGE"
$ws.Range("A38").Value = "RETRI PALADIN"

# Game 6 (rows 39-42): Beast Hunter / Demon Lock / Holy Priest -> WIN 2 / WIN 3
$ws.Range("A39").Value = 6
$ws.Range("A40").Value = "BEAST HUNTER"
$ws.Range("C40").Value = "WIN 2"
$ws.Range("A41").Value = "DEMON LOCK"
$ws.Range("C41").Value = "WIN 3"
$ws.Range("A42").Value = "HOLY PRIEST"

$ws.Range("D29").Select()
